$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency values in column C
$ws.Range("C2").Value = 3861
$ws.Range("C3").Value = 3664
$ws.Range("C4").Value = 2816
$ws.Range("C5").Value = 2003
$ws.Range("C6").Value = 1743
$ws.Range("C7").Value = 799
$ws.Range("C8").Value = 592
$ws.Range("C9").Value = 542
$ws.Range("C10").Value = 528
$ws.Range("C11").Value = 517

# Update Complementary Product category for row 11
$ws.Range("B11").Value = "Seasonal & Holidays"
